$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - column F is "想去人数" (number of people wanting to go)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 132
$ws1.Range("F9").Value = 2310
$ws1.Range("F10").Value = 119
$ws1.Range("F11").Value = 65
$ws1.Range("F13").Value = 1409
$ws1.Range("F14").Value = 499
$ws1.Range("F17").Value = 216
$ws1.Range("F26").Value = 1438
$ws1.Range("F29").Value = 176
$ws1.Range("F30").Value = 181

# Sheet "全部类型" (all types) - same underlying rows, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 133
$ws4.Range("F10").Value = 2310
$ws4.Range("F11").Value = 119
$ws4.Range("F12").Value = 65
$ws4.Range("F14").Value = 1409
$ws4.Range("F15").Value = 499
$ws4.Range("F18").Value = 216
$ws4.Range("F27").Value = 1438
$ws4.Range("F30").Value = 176
$ws4.Range("F31").Value = 181
